$wb = $excel.ActiveWorkbook

function Remove-CellHyperlink($ws, $addr) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.Delete()
        }
    }
}

# --- Overview sheet: update status text for the handed-off file row ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Handoff failed"
$wsOverview.Range("C2").Value = "Handoff failed"

# --- zh-cn sheet: clear the failed handoff's latest-handoff link/date, mark reason Ignored ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B2").Value = "Handoff failed"
Remove-CellHyperlink $wsZhCn '$C$2'
$wsZhCn.Range("C2").Clear()
$wsZhCn.Range("D2").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("H2").Value = "Ignored"

# --- de-de sheet: same treatment ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B2").Value = "Handoff failed"
Remove-CellHyperlink $wsDeDe '$C$2'
$wsDeDe.Range("C2").Clear()
$wsDeDe.Range("D2").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("H2").Value = "Ignored"
